$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("caña")

$desc = "• Medida: 3.90m en 4 secciones`n• Composición: grafito IM6`n• 5 pasahilos SiC de titanio de 3 y 4 patas, aptos para multifilamento`n• Acción: Heavy`n• Capacidad de lanzamiento máximo: 100 a 220 gramos`n• Peso: 480 gramos`n• Medida cerrada: 122 cm`n• Diámetro: Punta 3.2 mm / Base 26 mm`n• Distancia entre el taco y el porta reel: 65 cm`n• Primer pasahilo rebatible para transporte`n• Empuñadura de polímero anti-deslizante`n• Exclusivo tope interno de goma que evita el impacto entre sus tramos`n• Portareel de grafito a rosca con capuchones de acero inoxidable`n• Incluye capuchón protector de pasahilos y funda de tela"

$ws.Range("C2").Value2 = $desc
$ws.Range("C2").WrapText = $true

$ws.Columns.Item(3).ColumnWidth = 12.85546875
$ws.Columns.Item(5).ColumnWidth = 82.85546875

$ws.Rows.Item(2).RowHeight = 409.6

$ws.Range("E27").Select()
